# Update odds data for 2024-10-14 FlashScore weekly games workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (Once Caldas vs Dep. Pasto) ---
$ws.Range("G5").Value  = 1.91
$ws.Range("H5").Value  = 3.3
$ws.Range("I5").Value  = 4.2
$ws.Range("J5").Value  = 2.75
$ws.Range("L5").Value  = 5
$ws.Range("M5").Value  = 1.11
$ws.Range("N5").Value  = 6.5
$ws.Range("U5").Value  = 2.25
$ws.Range("V5").Value  = 1.57
$ws.Range("W5").Value  = 5.5
$ws.Range("AE5").Value = 21
$ws.Range("AH5").Value = 19
$ws.Range("AN5").Value = 3.75
$ws.Range("AP5").Value = 26
$ws.Range("AS5").Value = 251
$ws.Range("AX5").Value = 26

# --- Row 6 (America De Cali vs Ind. Medellin) ---
$ws.Range("I6").Value  = 4.75
$ws.Range("J6").Value  = 2.6
$ws.Range("K6").Value  = 2
$ws.Range("U6").Value  = 2.2
$ws.Range("V6").Value  = 1.62
$ws.Range("W6").Value  = 5.5
$ws.Range("X6").Value  = 7.5
$ws.Range("Z6").Value  = 15
$ws.Range("AC6").Value = 7
$ws.Range("AH6").Value = 21
$ws.Range("AZ6").Value = 101
